$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: for numeric-looking "Price" text values, force the cell
# to Text format before writing so Excel keeps the literal string (e.g. "1.00",
# "0.330", "71.004.80") instead of silently re-parsing it as a number, then clear
# the temporary formatting so the cell keeps the workbook default (unstyled) look.

# Rows 42 and 43: coin listing order changed (USDe <-> dogwifhat) with updated data
$ws.Cells.Item(42, 2).Value = "dogwifhat"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$dCell42 = $ws.Cells.Item(42, 4)
$dCell42.NumberFormat = "@"
$dCell42.Value = "2.63"
$dCell42.ClearFormats()
$ws.Cells.Item(42, 5).Value = "  +9.53%  "
$ws.Cells.Item(43, 2).Value = "USDe"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$dCell43 = $ws.Cells.Item(43, 4)
$dCell43.NumberFormat = "@"
$dCell43.Value = "1.00"
$dCell43.ClearFormats()
$ws.Cells.Item(43, 5).Value = "  +0.02%  "

# Remaining rows: refreshed Price (column D)
$priceUpdates = @{
    2 = "71.004.80"
    3 = "2.614.29"
    5 = "604.23"
    6 = "180.11"
    8 = "0.525"
    9 = "2.613.54"
    10 = "0.165"
    12 = "0.348"
    15 = "26.55"
    17 = "70.986.21"
    18 = "2.608.04"
    19 = "380.69"
    20 = "11.52"
    21 = "7.76"
    23 = "72.36"
    24 = "4.45"
    27 = "9.57"
    28 = "2.712.30"
    29 = "1.00"
    30 = "0.0₃0950"
    31 = "529.83"
    32 = "8.02"
    35 = "1.00"
    36 = "164.01"
    39 = "1.88"
    40 = "18.95"
    44 = "5.03"
    45 = "0.330"
    46 = "40.05"
    47 = "154.48"
    48 = "3.64"
    50 = "1.68"
    51 = "0.0₆0263"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.ClearFormats()
}

# Remaining rows: refreshed Volume(1h) (column E)
$volumeUpdates = @{
    2 = "  +3.77%  "
    3 = "  +4.05%  "
    4 = "  +0.04%  "
    5 = "  +2.05%  "
    6 = "  +2.20%  "
    7 = "  -0.04%  "
    8 = "  +1.63%  "
    9 = "  +4.03%  "
    10 = "  +14.67%  "
    11 = "  +0.13%  "
    12 = "  +3.05%  "
    13 = "  +0.14%  "
    15 = "  +3.00%  "
    16 = "  +7.18%  "
    17 = "  +4.09%  "
    18 = "  +4.06%  "
    19 = "  +8.62%  "
    20 = "  +5.01%  "
    21 = "  +3.42%  "
    22 = "  -0.73%  "
    23 = "  +1.69%  "
    24 = "  +5.50%  "
    25 = "  +0.07%  "
    26 = "  +6.27%  "
    27 = "  +5.14%  "
    28 = "  +4.23%  "
    29 = "  +0.07%  "
    30 = "  +5.50%  "
    31 = "  +3.96%  "
    32 = "  +2.54%  "
    33 = "  +3.32%  "
    34 = "  +2.64%  "
    35 = "  -0.01%  "
    36 = "  +1.33%  "
    37 = "  -1.79%  "
    38 = "  +4.28%  "
    39 = "  +6.75%  "
    44 = "  +4.11%  "
    45 = "  +0.45%  "
    46 = "  +2.68%  "
    47 = "  +2.69%  "
    48 = "  +2.14%  "
    49 = "  +2.39%  "
    50 = "  +4.91%  "
    51 = "  +1.38%  "
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}
